$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Range("C6").Value = 222.0942
$ws.Range("C9").Value = 278.08589999999998
$ws.Range("C13").Value = 10.607900000000001
$ws.Range("D13").Value = 10.607900000000001
$ws.Range("H13").Value = 5.7026000000000003
$ws.Range("I13").Value = 3.0617999999999999
$ws.Range("H20").Formula = "=C6/H6"
$ws.Range("H23").Formula = "=C9/H9"
$ws.Range("H24").Formula = "=C10/H10"
$ws.Range("H27").Formula = "=C13/H13"
